# ConsultaDetalleFactura.xlsx - "Actualizacion data y funcionalidades"
#
# Updates the "Datos" sheet: refreshes the cached session/transaction id
# (column N, rows 2-3) used by the data-driven test, widens column D to fit
# its content, and leaves the active selection on L7 (matching the state
# the workbook was saved in).

$wb = $excel.ActiveWorkbook

$wsDatos = $wb.Worksheets.Item("Datos")
$wsDatos.Activate()

# --- Data update -----------------------------------------------------
# N2 / N3 held a stale numeric id (63337); bump it to the current value.
$wsDatos.Range("N2").Value = 65468
$wsDatos.Range("N3").Value = 65468

# --- Column width ------------------------------------------------------
# Column D ("usuario" / "recaudosnatik66") did not have an explicit width
# before; size it to fit its contents (best-fit ~14.55 chars).
$wsDatos.Columns.Item(4).ColumnWidth = 13.65

# --- Selection ----------------------------------------------------------
# Workbook was left scrolled back to the top-left with L7 selected
# (previously the view was pinned to D1 with N2 selected).
$wsDatos.Range("L7").Select()
